$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update OT (column E) codes that changed from "Pendiente ADM"/ICD number to real OT numbers ---
$ws.Cells.Item(2, 5).Value = '''02281899 '   # Caso 6262
$ws.Cells.Item(5, 5).Value = '''02289539 '   # Caso 8004
$ws.Cells.Item(7, 5).Value = '''02271594 '   # Caso 8029
$ws.Cells.Item(8, 5).Value = '''02289689 '   # Caso 8048

# --- Remove rows that were taken off the board (Z7, Z2, Z5, Z8), bottom-to-top so row numbers stay valid ---
$ws.Rows(16).Delete() # Z8 / ALVAREZ JONTE AV. 1808
$ws.Rows(14).Delete() # Z5 / SERRANO 1074
$ws.Rows(11).Delete() # Z2 / LAVALLEJA 1030
$ws.Rows(9).Delete()  # Z7 / JULIAN ALVAREZ 928

# --- After the deletions, "S01335742" (Tinogasta 5182) landed on row 14; update its OT code too ---
$ws.Cells.Item(14, 5).Value = '''02289819 '   # Caso S01335742

# --- Append the new rows reported for the 1/14 and 1/16 batches ---
# Row 15: Caso 8133
$ws.Cells.Item(15, 1).Value = '''8133'
$ws.Cells.Item(15, 2).Value = '''1/14/2026'
$ws.Cells.Item(15, 3).Value = 'YERBAL 5612'
$ws.Cells.Item(15, 4).Value = 10
$ws.Cells.Item(15, 5).Value = '''02282092 '
$ws.Cells.Item(15, 6).Value = 'Optical Power'
$ws.Cells.Item(15, 7).Value = 'Pendiente'
$ws.Cells.Item(15, 8).Value = 'tendido bajo'
$ws.Cells.Item(15, 9).Value = 1
$ws.Cells.Item(15, 10).Value = '{"direccionesNormalizadas": [{"altura": 5612, "cod_calle": 26003, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.505213", "y": "-34.637279"}, "direccion": "YERBAL 5612, CABA", "nombre_calle": "YERBAL", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(15, 11).Value = -58.505213
$ws.Cells.Item(15, 12).Value = -34.637279
$ws.Cells.Item(15, 13).Value = 'Devoto'
$ws.Cells.Item(15, 14).Value = 'Capital Norte'

# Row 16: Caso 8134
$ws.Cells.Item(16, 1).Value = '''8134'
$ws.Cells.Item(16, 2).Value = '''1/14/2026'
$ws.Cells.Item(16, 3).Value = 'ESTADO DE PALESTINA 511'
$ws.Cells.Item(16, 4).Value = 5
$ws.Cells.Item(16, 5).Value = '''Pendiente ADM'
$ws.Cells.Item(16, 6).Value = 'Optical Power'
$ws.Cells.Item(16, 7).Value = 'Pendiente'
$ws.Cells.Item(16, 8).Value = 'tendido bajo'
$ws.Cells.Item(16, 9).Value = 1
$ws.Cells.Item(16, 10).Value = '{"direccionesNormalizadas": [{"altura": 511, "cod_calle": 19016, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.426036", "y": "-34.605725"}, "direccion": "ESTADO DE PALESTINA 511, CABA", "nombre_calle": "ESTADO DE PALESTINA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(16, 11).Value = -58.426036
$ws.Cells.Item(16, 12).Value = -34.605725
$ws.Cells.Item(16, 13).Value = 'Almagro'
$ws.Cells.Item(16, 14).Value = 'Capital Sur'

# Row 17: Caso 8135
$ws.Cells.Item(17, 1).Value = '''8135'
$ws.Cells.Item(17, 2).Value = '''1/14/2026'
$ws.Cells.Item(17, 3).Value = 'PERON, JUAN DOMINGO, TTE. GENERAL 4010'
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = '''Pendiente ADM'
$ws.Cells.Item(17, 6).Value = 'Optical Power'
$ws.Cells.Item(17, 7).Value = 'Pendiente'
$ws.Cells.Item(17, 8).Value = 'tendido bajo'
$ws.Cells.Item(17, 9).Value = 1
$ws.Cells.Item(17, 10).Value = '{"direccionesNormalizadas": [{"altura": 4010, "cod_calle": 3050, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.422975", "y": "-34.605999"}, "direccion": "PERON, JUAN DOMINGO, TTE. GENERAL 4010, CABA", "nombre_calle": "PERON, JUAN DOMINGO, TTE. GENERAL", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(17, 11).Value = -58.422975
$ws.Cells.Item(17, 12).Value = -34.605999
$ws.Cells.Item(17, 13).Value = 'Almagro'
$ws.Cells.Item(17, 14).Value = 'Capital Sur'

# Row 18: Caso S00922329
$ws.Cells.Item(18, 1).Value = '''S00922329'
$ws.Cells.Item(18, 2).Value = '''1/16/2026'
$ws.Cells.Item(18, 3).Value = 'PUEYRREDON, HONORIO, DR. AV. 620'
$ws.Cells.Item(18, 4).Value = 6
$ws.Cells.Item(18, 5).Value = '''Pendiente ADM'
$ws.Cells.Item(18, 6).Value = 'Optical Power'
$ws.Cells.Item(18, 7).Value = 'Pendiente'
$ws.Cells.Item(18, 8).Value = 'tendido bajo'
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(18, 10).Value = '{"direccionesNormalizadas": [{"altura": 620, "cod_calle": 17133, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.441834", "y": "-34.612760"}, "direccion": "PUEYRREDON, HONORIO, DR. AV. 620, CABA", "nombre_calle": "PUEYRREDON, HONORIO, DR. AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(18, 11).Value = -58.441834
$ws.Cells.Item(18, 12).Value = -34.61276
$ws.Cells.Item(18, 13).Value = 'Almagro'
$ws.Cells.Item(18, 14).Value = 'Capital Sur'

# Row 19: Caso S00984490
$ws.Cells.Item(19, 1).Value = '''S00984490'
$ws.Cells.Item(19, 2).Value = '''1/16/2026'
$ws.Cells.Item(19, 3).Value = 'HABANA 2235'
$ws.Cells.Item(19, 4).Value = 12
$ws.Cells.Item(19, 5).Value = '''Pendiente ADM'
$ws.Cells.Item(19, 6).Value = 'Optical Power'
$ws.Cells.Item(19, 7).Value = 'Pendiente'
$ws.Cells.Item(19, 8).Value = 'cable en panza'
$ws.Cells.Item(19, 9).Value = 1
$ws.Cells.Item(19, 10).Value = '{"direccionesNormalizadas": [{"altura": 2235, "cod_calle": 8001, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.491963", "y": "-34.584559"}, "direccion": "HABANA 2235, CABA", "nombre_calle": "HABANA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(19, 11).Value = -58.491963
$ws.Cells.Item(19, 12).Value = -34.584559
$ws.Cells.Item(19, 13).Value = 'Paternal'
$ws.Cells.Item(19, 14).Value = 'Capital Norte'

# Row 20: Caso S01108235
$ws.Cells.Item(20, 1).Value = '''S01108235'
$ws.Cells.Item(20, 2).Value = '''1/16/2026'
$ws.Cells.Item(20, 3).Value = 'PUEYRREDON, HONORIO, DR. AV. 632'
$ws.Cells.Item(20, 4).Value = 6
$ws.Cells.Item(20, 5).Value = '''Pendiente ADM'
$ws.Cells.Item(20, 6).Value = 'Optical Power'
$ws.Cells.Item(20, 7).Value = 'Pendiente'
$ws.Cells.Item(20, 8).Value = 'tendido bajo'
$ws.Cells.Item(20, 9).Value = 1
$ws.Cells.Item(20, 10).Value = '{"direccionesNormalizadas": [{"altura": 632, "cod_calle": 17133, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.441881", "y": "-34.612693"}, "direccion": "PUEYRREDON, HONORIO, DR. AV. 632, CABA", "nombre_calle": "PUEYRREDON, HONORIO, DR. AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(20, 11).Value = -58.441881
$ws.Cells.Item(20, 12).Value = -34.612693
$ws.Cells.Item(20, 13).Value = 'Almagro'
$ws.Cells.Item(20, 14).Value = 'Capital Sur'

# Row 21: Caso S01335725
$ws.Cells.Item(21, 1).Value = '''S01335725'
$ws.Cells.Item(21, 2).Value = '''1/16/2026'
$ws.Cells.Item(21, 3).Value = 'TINOGASTA 5196'
$ws.Cells.Item(21, 4).Value = 11
$ws.Cells.Item(21, 5).Value = '''Pendiente ADM'
$ws.Cells.Item(21, 6).Value = 'Optical Power'
$ws.Cells.Item(21, 7).Value = 'Pendiente'
$ws.Cells.Item(21, 8).Value = 'tendido bajo'
$ws.Cells.Item(21, 9).Value = 1
$ws.Cells.Item(21, 10).Value = '{"direccionesNormalizadas": [{"altura": 5196, "cod_calle": 21032, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.519645", "y": "-34.615857"}, "direccion": "TINOGASTA 5196, CABA", "nombre_calle": "TINOGASTA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Cells.Item(21, 11).Value = -58.519645
$ws.Cells.Item(21, 12).Value = -34.615857
$ws.Cells.Item(21, 13).Value = 'Devoto'
$ws.Cells.Item(21, 14).Value = 'Capital Norte'

